$p = $ppt.ActivePresentation

# --- 1. Update the table style on all 13 "Convolutional Neural Networks" slides ---
$tableSlideNumbers = @(15,16,17,18,19,21,22,25,26,27,28,31,33)
foreach ($n in $tableSlideNumbers) {
    $slide = $p.Slides.Item($n)
    $tableShape = $slide.Shapes.Item(4)
    $tableShape.Table.ApplyStyle("{988620D0-7669-41F4-938F-0C0862509A25}")
}

# --- 2. Slide 34: fix the lab title, logo placement, and notebook link ---
$slide34 = $p.Slides.Item(34)

# Title text: "Code Lab #2" -> "Lab Exercise #2"
$titleShape = $slide34.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Convolutional Neural Networks - Lab Exercise #2"

# Logo picture: reposition/resize
$logoShape = $slide34.Shapes.Item(2)
$logoShape.Left = 0
$logoShape.Top = 4.084645669291339
$logoShape.Width = 93.23228346456693
$logoShape.Height = 46.452757

# Notebook hyperlink: merge the two runs into one correctly-labeled hyperlink run
$bodyShape = $slide34.Shapes.Item(3)
$bodyTextRange = $bodyShape.TextFrame.TextRange
$linkRange = $bodyTextRange.Characters(77, 53)
$linkRange.Text = "Deep Learning by Design - Workshop - Chapter 2.ipynb"
$newLinkRange = $bodyTextRange.Characters(77, 54)
$newLinkRange.Font.Underline = -1
$newLinkRange.Font.Color.ObjectThemeColor = 11
$hyperlink = $newLinkRange.ActionSettings(1).Hyperlink
$hyperlink.Address = "https://github.com/GoogleCloudPlatform/keras-idiomatic-programmer/blob/master/workshops/Idiomatic%20Programmer%20-%20handbook%201%20-%20Codelab%202.ipynb"
